# Adding My Data Inside Excel File
# Fills row 2 with the submitter's name, email and repo link, turning the
# email and repo link into hyperlinks (mirrors the author's manual entry in
# Excel, which auto-applies the built-in "Hyperlink" cell style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data -----------------------------------------------------------
$ws.Range("A2").Value = "أسامة محمد عبدالمنعم على"
$ws.Range("B2").Value = "asamaaly70@gmail.com"
$ws.Range("C2").Value = "https://github.com/WalTeR-RE/Security-Task.git"

# --- Hyperlinks (email + repo link) ---------------------------------------
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:asamaaly70@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/WalTeR-RE/Security-Task.git")

# --- Column widths, widened to fit the new content -------------------------
$ws.Columns("A:A").ColumnWidth = 25
$ws.Columns("B:B").ColumnWidth = 30.7109375
$ws.Columns("C:C").ColumnWidth = 46.5703125

# --- Selection, left where the author last clicked -------------------------
[void]$ws.Range("C8").Select()
